$d = $word.ActiveDocument

$replacements = @(
    @("90-59=", "54+24="),
    @("95-72=", "18+37="),
    @("90-34=", "28-12="),
    @("6+52=", "45-7="),
    @("85-1=", "77-55="),
    @("68+0=", "98-92="),
    @("69-42=", "32-20="),
    @("60+7=", "8+29="),
    @("13+77=", "48+31="),
    @("75+10=", "59+0="),
    @("45-21=", "75-13="),
    @("11+9=", "59-38="),
    @("65+33=", "99-88="),
    @("98-26=", "13+44="),
    @("18-5=", "60+8="),
    @("97-45=", "68-63="),
    @("68-14=", "92-30="),
    @("77-72=", "37-21="),
    @("52-42=", "10-3="),
    @("68+6=", "8+36="),
    @("40+1=", "93-17="),
    @("54-46=", "13+84="),
    @("25+6=", "35+22="),
    @("64-9=", "77-52="),
    @("71+25=", "47-18="),
    @("52+45=", "69-29="),
    @("30+31=", "32+45="),
    @("10-5=", "35+62="),
    @("7+15=", "47-28="),
    @("85-75=", "72-67="),
    @("86-64=", "78-21="),
    @("91-10=", "87-51="),
    @("96-91=", "49+46="),
    @("36-2=", "75-46="),
    @("77-40=", "4+78="),
    @("16+50=", "47-25="),
    @("17+61=", "49-31="),
    @("72-63=", "25+60="),
    @("75-67=", "1-0="),
    @("98-22=", "21-18="),
    @("43+38=", "47-40="),
    @("25+47=", "42+52="),
    @("7+58=", "6+42="),
    @("85+5=", "81-59="),
    @("44-44=", "55-48="),
    @("5+74=", "29+29="),
    @("96-25=", "5+21="),
    @("92-66=", "49-12="),
    @("42+10=", "73-23="),
    @("78-53=", "70-54="),
    @("76-53=", "82+0="),
    @("23+59=", "89-27="),
    @("93-38=", "82+4="),
    @("6+91=", "96-41="),
    @("27+4=", "95-39="),
    @("88-9=", "86-76="),
    @("4+77=", "13+25="),
    @("47-1=", "8+85="),
    @("83-18=", "84-55="),
    @("85-23=", "79-49="),
    @("72-19=", "96-18="),
    @("93-49=", "10+85="),
    @("61+3=", "1+80="),
    @("47-31=", "98-36="),
    @("35+37=", "38+6="),
    @("73+19=", "72+17="),
    @("90-64=", "68+16="),
    @("66-4=", "89+2="),
    @("91-36=", "72+0="),
    @("56-24=", "84-15="),
    @("27-24=", "41+39="),
    @("4+57=", "67-18="),
    @("40+22=", "20+27="),
    @("76-40=", "93-52="),
    @("27+6=", "90-76="),
    @("85-37=", "66+4="),
    @("95-36=", "90-48="),
    @("1+20=", "60+1="),
    @("84-58=", "87-1="),
    @("45-30=", "47+36="),
    @("16+37=", "22-18="),
    @("30-15=", "6+75="),
    @("82+12=", "39+0="),
    @("88-4=", "36+27="),
    @("13+56=", "64-46="),
    @("60-57=", "87-40="),
    @("21+25=", "55-19="),
    @("19+79=", "18+59="),
    @("63-59=", "83+7="),
    @("61-27=", "99-90="),
    @("17+35=", "83-19="),
    @("1+16=", "95-61="),
    @("55-30=", "32-10="),
    @("23+51=", "30+32="),
    @("15+38=", "52+6="),
    @("13+68=", "74-64="),
    @("90-47=", "33-24="),
    @("1+4=", "70-19="),
    @("94-72=", "94-7="),
    @("52+10=", "2+12=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
